$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("Steps to Reproduce if Requried")
$ws.Range("E1").EntireColumn.Delete()

# Delete columns L:M ("Testers comments", "Developer comments") which, after the
# previous delete, are the two trailing columns following the shift
$ws.Range("L1:M1").EntireColumn.Delete()

# Update the active cell selection to match the final saved state
$ws.Range("D9").Select()
